# Updated dataset uncertainty (correction to the scaling), implemented more script functionality
#
# The E column holds the "amp_unc" (current uncertainty) values, computed
# from the measured current in column C. The scaling coefficient/offset
# used in that computation was corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- E3: single cell, offset corrected from 0.00000005 to 0.00005 ---
$ws.Range("E3").Formula = "=0.0002*ABS(C3)+0.00005"

# --- E4:E14: shared formula block, same offset correction ---
$ws.Range("E4:E14").Formula = "=0.0002*ABS(C4)+0.00005"

# --- E15: now its own (unshared) formula; coefficient corrected 0.0002 -> 0.002 ---
$ws.Range("E15").Formula = "=0.002*ABS(C15)+0.00005"

# --- E16:E18: new shared formula block, coefficient corrected 0.0002 -> 0.002 ---
$ws.Range("E16:E18").Formula = "=0.002*ABS(C16)+0.00005"

# --- E19:E22: each its own formula, coefficient corrected 0.0002 -> 0.002,
#     offsets scaled up by 10x to match the new coefficient scale ---
$ws.Range("E19").Formula = "=0.002*ABS(C19)+0.0005"
$ws.Range("E20").Formula = "=0.002*ABS(C20)+0.005"
$ws.Range("E21").Formula = "=0.002*ABS(C21)+0.005"
$ws.Range("E22").Formula = "=0.002*ABS(C22)+0.05"

# --- E23:E27: new shared formula block, coefficient corrected 0.0002 -> 0.002 ---
$ws.Range("E23:E27").Formula = "=0.002*ABS(C23)+0.05"

# --- Update the active selection to reflect where the author ended up ---
$ws.Range("E27").Select()
